$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph (bold, themed blue, sz 32) holding a manual page
#    break, right after the "Base estrategica ..." paragraph and right
#    before the "Cursos e Certificacoes" heading paragraph.
# ---------------------------------------------------------------------------
$baseIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Base estrat") {
        $baseIdx = $i
    }
}

$baseP = $d.Paragraphs.Item($baseIdx)
$insertPos = $baseP.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$pageBreakParaXml = "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val=`"215E99`" w:themeColor=`"text2`" w:themeTint=`"BF`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val=`"215E99`" w:themeColor=`"text2`" w:themeTint=`"BF`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:br w:type=`"page`"/></w:r></w:p>"

$insertRange.InsertXML($pageBreakParaXml)

# ---------------------------------------------------------------------------
# 2) Move the rendered-page-break marker: add it to the "Cursos e
#    Certificacoes" heading run (it now starts the new page) and drop it
#    from the "Desenvolvedor Full Stack ..." run (no longer a page start).
# ---------------------------------------------------------------------------
$cursosIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Cursos e Certif") {
        $cursosIdx = $i
    }
}

$cursosP = $d.Paragraphs.Item($cursosIdx)
$cursosRange = $cursosP.Range
$cursosRange.Collapse(1)

$cursosXml = "<w:p $wNs $w14Ns w14:paraId=`"34669BBC`" w14:textId=`"173D8B8F`" w:rsidR=`"00C249CF`" w:rsidRDefault=`"00C249CF`" w:rsidP=`"00C249CF`"><w:pPr><w:pBdr><w:bottom w:val=`"single`" w:sz=`"6`" w:space=`"1`" w:color=`"auto`"/></w:pBdr><w:spacing w:before=`"240`" w:after=`"0`" w:line=`"276`" w:lineRule=`"auto`"/><w:rPr><w:b/><w:bCs/><w:color w:val=`"215E99`" w:themeColor=`"text2`" w:themeTint=`"BF`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr></w:pPr><w:r w:rsidRPr=`"00C249CF`"><w:rPr><w:b/><w:bCs/><w:color w:val=`"215E99`" w:themeColor=`"text2`" w:themeTint=`"BF`"/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Cursos e Certifica" + [char]0x00E7 + [char]0x00F5 + "es</w:t></w:r></w:p>"

$cursosRange.InsertXML($cursosXml)

$fullStackIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Desenvolvedor Full Stack") {
        $fullStackIdx = $i
    }
}

$fullStackP = $d.Paragraphs.Item($fullStackIdx)
$fullStackRange = $fullStackP.Range
$fullStackRange.Collapse(1)

$fullStackXml = "<w:p $wNs $w14Ns w14:paraId=`"5F713EE1`" w14:textId=`"77777777`" w:rsidR=`"00C249CF`" w:rsidRPr=`"00C249CF`" w:rsidRDefault=`"00C249CF`" w:rsidP=`"00C249CF`"><w:pPr><w:spacing w:before=`"240`" w:after=`"0`" w:line=`"276`" w:lineRule=`"auto`"/><w:rPr><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r w:rsidRPr=`"00C249CF`"><w:rPr><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Desenvolvedor Full Stack </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r w:rsidRPr=`"00C249CF`"><w:rPr><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>JavaScript</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r w:rsidRPr=`"00C249CF`"><w:rPr><w:b/><w:bCs/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> " + [char]0x2022 + " EBAC</w:t></w:r></w:p>"

$fullStackRange.InsertXML($fullStackXml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
